$wb = $excel.ActiveWorkbook

# Rename sheets (case adjustments / translation)
$wb.Worksheets.Item("exiobase").Name = "Exiobase"
$wb.Worksheets.Item("german").Name = "Deutsch"
$wb.Worksheets.Item("english").Name = "Englisch"

# Update selection on the "Exiobase" sheet (formerly "exiobase")
$wsExio = $wb.Worksheets.Item("Exiobase")
$wsExio.Range("F34").Select()

# "map" sheet keeps its own selection (F21) but is no longer the active tab
$wsMap = $wb.Worksheets.Item("map")
$wsMap.Range("F21").Select()

# "Englisch" (formerly "english") sheet becomes the active tab with a new selection
$wsEnglisch = $wb.Worksheets.Item("Englisch")
$wsEnglisch.Activate()
$wsEnglisch.Range("F31").Select()
